# Applies numeric updates to the Tiamat_Profits leve-profit data across all class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 486270.06
$ws.Range("J17").Value = 486270.06
$ws.Range("L17").Value = 1458810.18
$ws.Range("N17").Value = -1459146.18
# Row 112
$ws.Range("H112").Value = 1944.1538
$ws.Range("I112").Value = 583.3333
$ws.Range("J112").Value = 2352.4
$ws.Range("K112").Value = 1749.9999
$ws.Range("L112").Value = 7057.200000000001
$ws.Range("M112").Value = -641.9999
$ws.Range("N112").Value = -9273.200000000001
# Row 116
$ws.Range("H116").Value = 6672.5
$ws.Range("I116").Value = 9822.143
$ws.Range("J116").Value = 4668.1816
$ws.Range("K116").Value = 9822.143
$ws.Range("L116").Value = 4668.1816
$ws.Range("M116").Value = -6380.143
$ws.Range("N116").Value = -11552.1816
# Row 137
$ws.Range("H137").Value = 3997.1428
$ws.Range("I137").Value = 929.7143
$ws.Range("J137").Value = 6042.095
$ws.Range("K137").Value = 2789.1429
$ws.Range("L137").Value = 18126.285
$ws.Range("M137").Value = -239.1428999999998
$ws.Range("N137").Value = -23226.285
# Row 138
$ws.Range("H138").Value = 1845.12
$ws.Range("I138").Value = 1061.4642
$ws.Range("J138").Value = 2149.875
$ws.Range("K138").Value = 3184.3926
$ws.Range("L138").Value = 6449.625
$ws.Range("M138").Value = 1955.6074
$ws.Range("N138").Value = -16729.625
# Row 141
$ws.Range("H141").Value = 2244.262
$ws.Range("I141").Value = 1806.9667
$ws.Range("K141").Value = 5420.9001
$ws.Range("M141").Value = -240.9000999999998

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 1347.5
$ws.Range("I102").Value = 963.3333
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 963.3333
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 658.6667
$ws.Range("N102").Value = -5744
# Row 132
$ws.Range("H132").Value = 290496.12
$ws.Range("I132").Value = 5056.4443
$ws.Range("J132").Value = 1253855
$ws.Range("K132").Value = 15169.3329
$ws.Range("L132").Value = 3761565
$ws.Range("M132").Value = -12639.3329
$ws.Range("N132").Value = -3766625

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1991402.5
$ws.Range("I105").Value = 3981505
$ws.Range("K105").Value = 3981505
$ws.Range("M105").Value = -3979758
# Row 134
$ws.Range("H134").Value = 22246282
$ws.Range("I134").Value = 1141.9615
$ws.Range("J134").Value = 52687000
$ws.Range("K134").Value = 3425.8845
$ws.Range("L134").Value = 158061000
$ws.Range("M134").Value = -890.8844999999997
$ws.Range("N134").Value = -158066070

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 31201.8
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 38752.25
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 38752.25
$ws.Range("M3").Value = -887
$ws.Range("N3").Value = -38978.25
# Row 31
$ws.Range("H31").Value = 12242.216
$ws.Range("I31").Value = 14249.685
$ws.Range("J31").Value = 10123.223
$ws.Range("K31").Value = 14249.685
$ws.Range("L31").Value = 10123.223
$ws.Range("M31").Value = -13954.685
$ws.Range("N31").Value = -10713.223
# Row 34
$ws.Range("H34").Value = 12242.216
$ws.Range("I34").Value = 14249.685
$ws.Range("J34").Value = 10123.223
$ws.Range("K34").Value = 14249.685
$ws.Range("L34").Value = 10123.223
$ws.Range("M34").Value = -14047.685
$ws.Range("N34").Value = -10527.223
# Row 99
$ws.Range("H99").Value = 3227.6365
$ws.Range("I99").Value = 3038.625
$ws.Range("K99").Value = 3038.625
$ws.Range("M99").Value = -1540.625
# Row 126
$ws.Range("H126").Value = 3227.6365
$ws.Range("I126").Value = 3038.625
$ws.Range("K126").Value = 9115.875
$ws.Range("M126").Value = -6645.875

$ws = $wb.Worksheets.Item("CUL")
# Row 103
$ws.Range("H103").Value = 2034.1666
$ws.Range("I103").Value = 616.6667
$ws.Range("J103").Value = 2506.6667
$ws.Range("K103").Value = 1850.0001
$ws.Range("L103").Value = 7520.000100000001
$ws.Range("M103").Value = -971.0001
$ws.Range("N103").Value = -9278.000100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("L4").Value = 0
# Row 70
$ws.Range("H70").Value = 5887229.5
$ws.Range("I70").Value = 8337950
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 8337950
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -8337680
$ws.Range("N70").Value = -6040
# Row 73
$ws.Range("H73").Value = 5887229.5
$ws.Range("I73").Value = 8337950
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 8337950
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -8337014
$ws.Range("N73").Value = -7372
# Row 122
$ws.Range("H122").Value = 1202.8334
$ws.Range("I122").Value = 1154.4286
$ws.Range("J122").Value = 1372.25
$ws.Range("K122").Value = 3463.2858
$ws.Range("L122").Value = 4116.75
$ws.Range("M122").Value = -1013.2858
$ws.Range("N122").Value = -9016.75
# Row 135
$ws.Range("H135").Value = 64294.453
$ws.Range("J135").Value = 64294.453
$ws.Range("L135").Value = 64294.453
$ws.Range("N135").Value = -74434.45300000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2895.2222
$ws.Range("I7").Value = 3192.9
$ws.Range("J7").Value = 2523.125
$ws.Range("K7").Value = 3192.9
$ws.Range("L7").Value = 2523.125
$ws.Range("M7").Value = -3080.9
$ws.Range("N7").Value = -2747.125
# Row 22
$ws.Range("H22").Value = 1097.25
$ws.Range("I22").Value = 270
$ws.Range("J22").Value = 1924.5
$ws.Range("K22").Value = 270
$ws.Range("L22").Value = 1924.5
$ws.Range("M22").Value = 25
$ws.Range("N22").Value = -2514.5
# Row 27
$ws.Range("H27").Value = 1097.25
$ws.Range("I27").Value = 270
$ws.Range("J27").Value = 1924.5
$ws.Range("K27").Value = 270
$ws.Range("L27").Value = 1924.5
$ws.Range("M27").Value = -163
$ws.Range("N27").Value = -2138.5
# Row 40
$ws.Range("H40").Value = 101962.4
$ws.Range("I40").Value = 1596.8
$ws.Range("K40").Value = 1596.8
$ws.Range("M40").Value = -1460.8
# Row 82
$ws.Range("H82").Value = 1516.0714
$ws.Range("I82").Value = 1146.9286
$ws.Range("J82").Value = 1885.2142
$ws.Range("K82").Value = 1146.9286
$ws.Range("L82").Value = 1885.2142
$ws.Range("M82").Value = -785.9286
$ws.Range("N82").Value = -2607.2142
# Row 85
$ws.Range("H85").Value = 1516.0714
$ws.Range("I85").Value = 1146.9286
$ws.Range("J85").Value = 1885.2142
$ws.Range("K85").Value = 1146.9286
$ws.Range("L85").Value = 1885.2142
$ws.Range("M85").Value = 101.0714
$ws.Range("N85").Value = -4381.2142
# Row 122
$ws.Range("H122").Value = 2301.756
$ws.Range("I122").Value = 2287.5625
$ws.Range("J122").Value = 2352.2222
$ws.Range("K122").Value = 6862.6875
$ws.Range("L122").Value = 7056.6666
$ws.Range("M122").Value = -4412.6875
$ws.Range("N122").Value = -11956.6666
# Row 126
$ws.Range("H126").Value = 2895.2222
$ws.Range("I126").Value = 3192.9
$ws.Range("J126").Value = 2523.125
$ws.Range("K126").Value = 9578.700000000001
$ws.Range("L126").Value = 7569.375
$ws.Range("M126").Value = -7108.700000000001
$ws.Range("N126").Value = -12509.375

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 76251.5
$ws.Range("J5").Value = 76251.5
$ws.Range("L5").Value = 76251.5
$ws.Range("N5").Value = -76475.5
# Row 96
$ws.Range("H96").Value = 1132.4231
$ws.Range("I96").Value = 997
$ws.Range("J96").Value = 1317.091
$ws.Range("K96").Value = 997
$ws.Range("L96").Value = 1317.091
$ws.Range("M96").Value = 376
$ws.Range("N96").Value = -4063.091
# Row 126
$ws.Range("H126").Value = 709.04346
$ws.Range("I126").Value = 646.53845
$ws.Range("J126").Value = 790.3
$ws.Range("K126").Value = 1939.61535
$ws.Range("L126").Value = 2370.9
$ws.Range("M126").Value = 530.38465
$ws.Range("N126").Value = -7310.9
